$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.395.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.140.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.75%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.138.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.474"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.62%  "

$ws.Range("E13").Value = "  +1.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.667.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.52%  "

$ws.Range("E16").Value = "  +2.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.382.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.139.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.719"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.83%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.30%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.63%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.75%  "

$ws.Range("E30").Value = "  +3.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.85%  "

$ws.Range("E32").Value = "  -0.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.80%  "

$ws.Range("E35").Value = "  +0.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0743"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.35%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "449.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.59%  "

$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0394"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("E42").Value = "  +0.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.881.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.262"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.114"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.72%  "
